$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.504.74'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +3.58%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.733.19'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9996'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.12'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +3.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.0000'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4795'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +3.47%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +3.47%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06226'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.25%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.733.36'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +4.04%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07121'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +2.73%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.74'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +5.71%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +6.98%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.532'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +4.20%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '76.95'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.21%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.000'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.11%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.501.12'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +3.54%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.000'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.16%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000006913'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +3.00%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.72'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +2.95%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.958.05'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +4.42%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.579'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +3.23%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.897'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +3.25%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.329'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.84%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '136.19'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.42%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.34'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +2.66%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.797'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +4.33%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.412'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +2.77%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '106.66'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +2.22%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.989'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.79%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.07992'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +4.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.713'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +3.11%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04542'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +4.59%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.616'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.51%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6369'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +4.76%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9921'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +5.40%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9363'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.17%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '110.09'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +2.22%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.982'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +7.89%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.415'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.14%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.006'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.64%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.01507'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +3.73%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.706'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +13.15%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3906'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +5.10%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.943'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +13.59%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1192'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +7.42%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.27%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.900'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +3.57%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '30.74'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.267'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +5.08%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3424'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +3.06%  '
